$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block (rows 10-12): update scores and pick up the header-row style ---
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)

$ws.Range("B10").Value = 13
$ws.Range("D10").Value = 15
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 52
$ws.Range("E12").Value = "52/112"

# --- Fill in column A with the student answer for the rows that now have one,
#     picking up the "correctStyle" formatting already used by B10:B12 ---
$ws.Range("B10").Copy()
$aCells = @(16,18,20,23,25,26,30,31,32,38,39,40)
foreach ($r in $aCells) {
    $ws.Range("A$r").PasteSpecial(-4122)
}

$ws.Range("A16").Value = "Option A"
$ws.Range("A18").Value = "Option B"
$ws.Range("A20").Value = "Option B"
$ws.Range("A23").Value = "Option D"
$ws.Range("A25").Value = "Option A"
$ws.Range("A26").Value = "Option C"
$ws.Range("A30").Value = "Option B"
$ws.Range("A31").Value = "Option D"
$ws.Range("A32").Value = "Option C"
$ws.Range("A38").Value = "Option A"
$ws.Range("A39").Value = "Option D"
$ws.Range("A40").Value = "Option D"

# D18 also gets a value now, with the same correctStyle formatting
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Option D"

# --- Remove the now-unused duplicate "Student Ans/Correct Ans" block in columns G:H ---
$ws.Range("G:H").Delete()

# --- Remove the stray duplicate D/E answer data below row 18 (keep the header row 15 and rows 16-18) ---
$ws.Range("D19:E40").Clear()
